# "Building out contact pg" — append a new (still-empty) slide to the deck
# so work can start on the contact page, matching the dark theme used by
# the rest of the deck.

$p = $ppt.ActivePresentation

# Add a new Blank-layout slide (ppLayoutBlank = 12) at the end of the deck.
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 12)

# Give it the same solid dark background (#212529) the other slides use,
# instead of inheriting the (light) master/theme background.
$s.FollowMasterBackground = $false
$s.Background.Fill.ForeColor.RGB = 2696481   # RGB(0x21,0x25,0x29) = #212529
